$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.2
$ws.Range("I4").Value = 3.8
$ws.Range("R4").Value = 2.5
$ws.Range("S4").Value = 1.5
$ws.Range("U4").Value = 8.5
$ws.Range("X4").Value = 23
$ws.Range("G30").Value = 3
$ws.Range("I30").Value = 2.3
$ws.Range("K30").Value = 7.7
$ws.Range("L30").Value = 1.3
$ws.Range("N30").Value = 1.9
$ws.Range("O30").Value = 1.85
$ws.Range("P30").Value = 1.42
$ws.Range("Q30").Value = 2.7
$ws.Range("S30").Value = 2.05
$ws.Range("T30").Value = 9.75
$ws.Range("U30").Value = 17.5
$ws.Range("V30").Value = 11.25
$ws.Range("W30").Value = 40
$ws.Range("X30").Value = 27
$ws.Range("Z30").Value = 7.7
$ws.Range("AE30").Value = 7.8
$ws.Range("AI30").Value = 20
$ws.Range("AJ30").Value = 30
$ws.Range("T34").Value = 7.6
$ws.Range("X34").Value = 35
$ws.Range("AB34").Value = 18.5
$ws.Range("AE34").Value = 6.5
$ws.Range("AI34").Value = 24
$ws.Range("G35").Value = 2.62
$ws.Range("H35").Value = 3.1
$ws.Range("I35").Value = 2.75
$ws.Range("J35").Value = 1.09
$ws.Range("K35").Value = 6.7
$ws.Range("L35").Value = 1.39
$ws.Range("M35").Value = 2.85
$ws.Range("N35").Value = 2.15
$ws.Range("O35").Value = 1.65
$ws.Range("P35").Value = 1.47
$ws.Range("Q35").Value = 2.57
$ws.Range("R35").Value = 1.83
$ws.Range("S35").Value = 1.87
$ws.Range("T35").Value = 7.4
$ws.Range("V35").Value = 10.5
$ws.Range("X35").Value = 26
$ws.Range("Y35").Value = 40
$ws.Range("Z35").Value = 6.7
$ws.Range("AA35").Value = 6.2
$ws.Range("AB35").Value = 15.5
$ws.Range("AC35").Value = 80
$ws.Range("AD35").Value = 800
$ws.Range("AE35").Value = 7.8
$ws.Range("AG35").Value = 10.75
$ws.Range("AI35").Value = 27
$ws.Range("AJ35").Value = 40
$ws.Range("L39").Value = 1.4
$ws.Range("M39").Value = 2.55
$ws.Range("N39").Value = 2.15
$ws.Range("O39").Value = 1.55
$ws.Range("Q39").Value = 2.32
$ws.Range("R39").Value = 1.98
$ws.Range("S39").Value = 1.65
$ws.Range("Y39").Value = 55
$ws.Range("Z39").Value = 7.8
$ws.Range("AB39").Value = 18
$ws.Range("AC39").Value = 110
$ws.Range("AE39").Value = 6
$ws.Range("AJ39").Value = 37
$ws.Range("G53").Value = 1.83
$ws.Range("T53").Value = 7.5
$ws.Range("U53").Value = 8.5
$ws.Range("X53").Value = 11.25
$ws.Range("Z53").Value = 12
$ws.Range("AB53").Value = 10.75
$ws.Range("AC53").Value = 37
$ws.Range("AE53").Value = 10.5
$ws.Range("AF53").Value = 17.5
$ws.Range("AI53").Value = 24
$ws.Range("AJ53").Value = 26
$ws.Range("L56").Value = 1.3
$ws.Range("M56").Value = 3.4
$ws.Range("N56").Value = 1.98
$ws.Range("O56").Value = 1.83
$ws.Range("G58").Value = 3.9
$ws.Range("I58").Value = 1.7
$ws.Range("N58").Value = 1.29
$ws.Range("O58").Value = 3
$ws.Range("R58").Value = 1.33
$ws.Range("T58").Value = 26
$ws.Range("U58").Value = 35
$ws.Range("V58").Value = 14.5
$ws.Range("W58").Value = 70
$ws.Range("X58").Value = 29
$ws.Range("Z58").Value = 27
$ws.Range("AA58").Value = 10.25
$ws.Range("AC58").Value = 26
$ws.Range("AD58").Value = 100
$ws.Range("AE58").Value = 15
$ws.Range("AF58").Value = 13.5
$ws.Range("AH58").Value = 17
$ws.Range("AI58").Value = 11.75
$ws.Range("T59").Value = 11.75
$ws.Range("U59").Value = 16
$ws.Range("AE59").Value = 11
$ws.Range("J68").Value = 1.11
$ws.Range("K68").Value = 6.5
$ws.Range("N69").Value = 2.1
$ws.Range("O69").Value = 1.7
$ws.Range("N77").Value = 1.98
$ws.Range("O77").Value = 1.83
$ws.Range("G85").Value = 2.15
$ws.Range("H85").Value = 3.3
$ws.Range("I85").Value = 3.4
$ws.Range("J85").Value = 1.04
$ws.Range("K85").Value = 12
$ws.Range("L85").Value = 1.22
$ws.Range("M85").Value = 4
$ws.Range("R85").Value = 1.62
$ws.Range("S85").Value = 2.2
$ws.Range("T85").Value = 9
$ws.Range("U85").Value = 11
$ws.Range("W85").Value = 19
$ws.Range("X85").Value = 17
$ws.Range("Z85").Value = 12
$ws.Range("AB85").Value = 12
$ws.Range("AF85").Value = 19
$ws.Range("AG85").Value = 12
$ws.Range("AH85").Value = 34
$ws.Range("AI85").Value = 26
$ws.Range("AJ85").Value = 29
$ws.Range("I95").Value = 3.4
$ws.Range("T95").Value = 6.5
$ws.Range("AB95").Value = 15
$ws.Range("AG95").Value = 13
$ws.Range("G98").Value = 1.27
$ws.Range("H98").Value = 6
$ws.Range("I98").Value = 9.5
$ws.Range("N98").Value = 1.57
$ws.Range("O98").Value = 2.35
$ws.Range("R98").Value = 2
$ws.Range("S98").Value = 1.73
$ws.Range("W98").Value = 8
$ws.Range("Z98").Value = 15
$ws.Range("AB98").Value = 23
$ws.Range("AF98").Value = 41
$ws.Range("AG98").Value = 26
$ws.Range("J99").Value = 1.03
$ws.Range("K99").Value = 17
$ws.Range("R99").Value = 1.62
$ws.Range("S99").Value = 2.2
$ws.Range("T99").Value = 19
$ws.Range("Z99").Value = 17
$ws.Range("AA99").Value = 8.5
$ws.Range("AJ99").Value = 21
$ws.Range("G102").Value = 1.57
$ws.Range("H102").Value = 4.2
$ws.Range("T102").Value = 8.5
$ws.Range("AA102").Value = 8
$ws.Range("G103").Value = 1.67
$ws.Range("J103").Value = 1.11
$ws.Range("K103").Value = 6.5
$ws.Range("L103").Value = 1.44
$ws.Range("M103").Value = 2.63
$ws.Range("AA103").Value = 7
$ws.Range("AB103").Value = 23
$ws.Range("N107").Value = 2
$ws.Range("K110").Value = 17
$ws.Range("L110").Value = 1.17
$ws.Range("M110").Value = 5
$ws.Range("N110").Value = 1.57
$ws.Range("O110").Value = 2.35
$ws.Range("G124").Value = 2.32
$ws.Range("H124").Value = 3.1
$ws.Range("I124").Value = 3
$ws.Range("L124").Value = 1.33
$ws.Range("M124").Value = 2.77
$ws.Range("N124").Value = 1.98
$ws.Range("O124").Value = 1.65
$ws.Range("P124").Value = 1.4
$ws.Range("Q124").Value = 2.52
$ws.Range("R124").Value = 1.75
$ws.Range("S124").Value = 1.87
$ws.Range("T124").Value = 7.5
$ws.Range("U124").Value = 11.25
$ws.Range("V124").Value = 9
$ws.Range("W124").Value = 24
$ws.Range("X124").Value = 19.5
$ws.Range("Y124").Value = 30
$ws.Range("Z124").Value = 8.5
$ws.Range("AA124").Value = 6
$ws.Range("AB124").Value = 14
$ws.Range("AC124").Value = 70
$ws.Range("AD124").Value = 600
$ws.Range("AE124").Value = 8.5
$ws.Range("AF124").Value = 15
$ws.Range("AG124").Value = 10.75
$ws.Range("AH124").Value = 37
$ws.Range("AI124").Value = 28
$ws.Range("AJ124").Value = 37
